$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'42.768.20"
$ws.Range('E2').Value = "'  +0.80%  "
$ws.Range('D3').Value = "'2.317.64"
$ws.Range('E3').Value = "'  +1.71%  "
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = "'  -0.10%  "
$ws.Range('D5').Value = "'302.13"
$ws.Range('E5').Value = "'  -0.29%  "
$ws.Range('D6').Value = "'96.18"
$ws.Range('E6').Value = "'  +1.38%  "
$ws.Range('D7').Value = "'0.508"
$ws.Range('E7').Value = "'  +1.26%  "
$ws.Range('E8').Value = "'  +0.00%  "
$ws.Range('D9').Value = "'0.493"
$ws.Range('E9').Value = "'  +0.20%  "
$ws.Range('D10').Value = "'34.58"
$ws.Range('E10').Value = "'  -1.01%  "
$ws.Range('D11').Value = "'19.28"
$ws.Range('E11').Value = "'  +7.08%  "
$ws.Range('D12').Value = "'0.0786"
$ws.Range('E12').Value = "'  +1.28%  "
$ws.Range('D13').Value = "'0.120"
$ws.Range('E13').Value = "'  +0.58%  "
$ws.Range('D14').Value = "'6.75"
$ws.Range('E14').Value = "'  +1.32%  "
$ws.Range('D15').Value = "'2.676.63"
$ws.Range('E15').Value = "'  +1.55%  "
$ws.Range('D16').Value = "'2.305.39"
$ws.Range('E16').Value = "'  +0.79%  "
$ws.Range('D17').Value = "'0.786"
$ws.Range('E17').Value = "'  +2.20%  "
$ws.Range('D18').Value = "'42.714.71"
$ws.Range('E18').Value = "'  +0.87%  "
$ws.Range('D19').Value = "'12.27"
$ws.Range('E19').Value = "'  -3.16%  "
$ws.Range('D20').Value = "'6.14"
$ws.Range('E20').Value = "'  +3.17%  "
$ws.Range('D21').Value = "'0.0₃0891"
$ws.Range('E21').Value = "'  +0.67%  "
$ws.Range('D22').Value = "'68.06"
$ws.Range('E22').Value = "'  +1.56%  "
$ws.Range('E23').Value = "'  +5.33%  "
$ws.Range('D24').Value = "'235.85"
$ws.Range('E24').Value = "'  +0.17%  "
$ws.Range('E25').Value = "'  +0.09%  "
$ws.Range('D26').Value = "'2.42"
$ws.Range('E26').Value = "'  +1.79%  "
$ws.Range('D27').Value = "'24.38"
$ws.Range('E27').Value = "'  -1.10%  "
$ws.Range('E28').Value = "'  -1.03%  "
$ws.Range('D29').Value = "'166.42"
$ws.Range('D30').Value = "'9.13"
$ws.Range('E30').Value = "'  +2.44%  "
$ws.Range('D31').Value = "'32.50"
$ws.Range('E31').Value = "'  +1.00%  "
$ws.Range('E32').Value = "'  -0.09%  "
$ws.Range('D33').Value = "'5.02"
$ws.Range('E33').Value = "'  +2.14%  "
$ws.Range('D34').Value = "'17.70"
$ws.Range('E34').Value = "'  +1.08%  "
$ws.Range('D35').Value = "'4.45"
$ws.Range('E35').Value = "'  +0.19%  "
$ws.Range('D36').Value = "'0.0702"
$ws.Range('E36').Value = "'  +3.34%  "
$ws.Range('B38').Value = "'ARBITRUM"
$ws.Range('C38').Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range('D38').Value = "'1.77"
$ws.Range('E38').Value = "'  +3.01%  "
$ws.Range('B39').Value = "'Kaspa"
$ws.Range('C39').Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range('D39').Value = "'0.0997"
$ws.Range('E39').Value = "'  -0.32%  "
$ws.Range('D40').Value = "'0.109"
$ws.Range('E40').Value = "'  +0.60%  "
$ws.Range('D41').Value = "'2.71"
$ws.Range('E41').Value = "'  +2.50%  "
$ws.Range('D42').Value = "'20.11"
$ws.Range('E42').Value = "'  +12.95%  "
$ws.Range('D43').Value = "'1.952.74"
$ws.Range('E43').Value = "'  -1.48%  "
$ws.Range('D44').Value = "'10.41"
$ws.Range('E44').Value = "'  +4.28%  "
$ws.Range('D45').Value = "'0.0279"
$ws.Range('E45').Value = "'  +1.76%  "
$ws.Range('E46').Value = "'  +2.63%  "
$ws.Range('D47').Value = "'2.75"
$ws.Range('E47').Value = "'  +0.95%  "
$ws.Range('D48').Value = "'2.545.05"
$ws.Range('E48').Value = "'  +1.65%  "
$ws.Range('B49').Value = "'MultiversX"
$ws.Range('C49').Value = "'https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range('D49').Value = "'53.41"
$ws.Range('E49').Value = "'  +0.56%  "
$ws.Range('B50').Value = "'HuobiToken"
$ws.Range('C50').Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range('D50').Value = "'2.79"
$ws.Range('E50').Value = "'  -4.68%  "
$ws.Range('D51').Value = "'72.13"
$ws.Range('E51').Value = "'  +2.61%  "
